# Update "想去人数" (interested count) values in sheet "展览" and "全部类型"
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    5  = 2817
    9  = 1508
    10 = 36
    11 = 73
    18 = 48
    20 = 83
    22 = 2790
    24 = 7
    25 = 51
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
